# chore: update Sheets via scheduled runner
# Refresh market-price / profit figures (columns H-N) on the Yojimbo_Profits
# workbook's per-job leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 599.10345
$ws.Range("I107").Value = 486.05264
$ws.Range("J107").Value = 813.9
$ws.Range("K107").Value = 486.05264
$ws.Range("L107").Value = 813.9
$ws.Range("M107").Value = 1433.94736
$ws.Range("N107").Value = -4653.9
$ws.Range("H120").Value = 149999.33
$ws.Range("J120").Value = 149999.33
$ws.Range("L120").Value = 149999.33
$ws.Range("N120").Value = -159675.33
$ws.Range("H125").Value = 1021.0909
$ws.Range("I125").Value = 832.6667
$ws.Range("J125").Value = 1247.2
$ws.Range("K125").Value = 7494.0003
$ws.Range("L125").Value = 11224.8
$ws.Range("M125").Value = -5034.0003
$ws.Range("N125").Value = -16144.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 28500
$ws.Range("J10").Value = 28500
$ws.Range("L10").Value = 28500
$ws.Range("N10").Value = -28840
$ws.Range("H74").Value = 888.5762999999999
$ws.Range("I74").Value = 594.73334
$ws.Range("J74").Value = 1833.0714
$ws.Range("K74").Value = 594.73334
$ws.Range("L74").Value = 1833.0714
$ws.Range("M74").Value = 279.26666
$ws.Range("N74").Value = -3581.0714
$ws.Range("H77").Value = 888.5762999999999
$ws.Range("I77").Value = 594.73334
$ws.Range("J77").Value = 1833.0714
$ws.Range("K77").Value = 2973.6667
$ws.Range("L77").Value = 9165.357
$ws.Range("M77").Value = 1394.3333
$ws.Range("N77").Value = -17901.357
$ws.Range("H122").Value = 7695143.5
$ws.Range("I122").Value = 12502246
$ws.Range("J122").Value = 3780
$ws.Range("K122").Value = 37506738
$ws.Range("L122").Value = 11340
$ws.Range("M122").Value = -37504288
$ws.Range("N122").Value = -16240
$ws.Range("H124").Value = 9031.166999999999
$ws.Range("J124").Value = 9031.166999999999
$ws.Range("L124").Value = 9031.166999999999
$ws.Range("N124").Value = -18851.167
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H129").Value = 46014.285
$ws.Range("J129").Value = 46014.285
$ws.Range("L129").Value = 46014.285
$ws.Range("N129").Value = -56014.285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 759.28
$ws.Range("I3").Value = 694.381
$ws.Range("J3").Value = 1100
$ws.Range("K3").Value = 694.381
$ws.Range("L3").Value = 1100
$ws.Range("M3").Value = -580.381
$ws.Range("N3").Value = -1328

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1828.5
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 2157
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 2157
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -5153
$ws.Range("H126").Value = 1828.5
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 2157
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 6471
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -11411
$ws.Range("H127").Value = 32430
$ws.Range("J127").Value = 32430
$ws.Range("L127").Value = 32430
$ws.Range("N127").Value = -42350
$ws.Range("H129").Value = 45199.75
$ws.Range("J129").Value = 45199.75
$ws.Range("L129").Value = 45199.75
$ws.Range("N129").Value = -55199.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 8808.916999999999
$ws.Range("J107").Value = 561.55554
$ws.Range("L107").Value = 1684.66662
$ws.Range("N107").Value = -5524.66662
$ws.Range("H110").Value = 4092.4546
$ws.Range("I110").Value = 1853.5
$ws.Range("J110").Value = 4590
$ws.Range("K110").Value = 5560.5
$ws.Range("L110").Value = 13770
$ws.Range("M110").Value = -1470.5
$ws.Range("N110").Value = -21950
$ws.Range("H113").Value = 1319.8235
$ws.Range("I113").Value = 1429.1333
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 4287.3999
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -2117.3999
$ws.Range("N113").Value = -5840
$ws.Range("H115").Value = 3156.0667
$ws.Range("J115").Value = 3310.0715
$ws.Range("L115").Value = 9930.2145
$ws.Range("N115").Value = -12280.2145
$ws.Range("H118").Value = 1920.2084
$ws.Range("I118").Value = 603.125
$ws.Range("K118").Value = 1809.375
$ws.Range("M118").Value = -566.375
$ws.Range("H119").Value = 4560.6665
$ws.Range("I119").Value = 5133.3335
$ws.Range("J119").Value = 3988
$ws.Range("K119").Value = 15400.0005
$ws.Range("L119").Value = 11964
$ws.Range("M119").Value = -10562.0005
$ws.Range("N119").Value = -21640
$ws.Range("H120").Value = 20723.572
$ws.Range("J120").Value = 20723.572
$ws.Range("L120").Value = 62170.716
$ws.Range("N120").Value = -71846.716
$ws.Range("H131").Value = 708.8163500000001
$ws.Range("J131").Value = 895.5
$ws.Range("L131").Value = 2686.5
$ws.Range("N131").Value = -12766.5
$ws.Range("H132").Value = 1119.3889
$ws.Range("I132").Value = 816
$ws.Range("J132").Value = 1422.7778
$ws.Range("K132").Value = 7344
$ws.Range("L132").Value = 12805.0002
$ws.Range("M132").Value = -4814
$ws.Range("N132").Value = -17865.0002
$ws.Range("H138").Value = 2032.5
$ws.Range("I138").Value = 871.4286
$ws.Range("J138").Value = 2460.2632
$ws.Range("K138").Value = 2614.2858
$ws.Range("L138").Value = 7380.7896
$ws.Range("M138").Value = 2525.7142
$ws.Range("N138").Value = -17660.7896
$ws.Range("H139").Value = 1112
$ws.Range("I139").Value = 641.86664
$ws.Range("J139").Value = 2875
$ws.Range("K139").Value = 1925.59992
$ws.Range("L139").Value = 8625
$ws.Range("M139").Value = 3214.40008
$ws.Range("N139").Value = -18905
$ws.Range("H140").Value = 1789.3256
$ws.Range("I140").Value = 1016.2
$ws.Range("K140").Value = 3048.6
$ws.Range("M140").Value = 2131.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 28500000
$ws.Range("I11").Value = 4666666.5
$ws.Range("J11").Value = 100000000
$ws.Range("K11").Value = 4666666.5
$ws.Range("L11").Value = 100000000
$ws.Range("M11").Value = -4666527.5
$ws.Range("N11").Value = -100000278
$ws.Range("H126").Value = 1480.1111
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1480.1111
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 4440.3333
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -9380.3333
$ws.Range("H129").Value = 43966.668
$ws.Range("J129").Value = 43966.668
$ws.Range("L129").Value = 43966.668
$ws.Range("N129").Value = -53966.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2783.8
$ws.Range("I7").Value = 1204
$ws.Range("J7").Value = 3178.75
$ws.Range("K7").Value = 1204
$ws.Range("L7").Value = 3178.75
$ws.Range("M7").Value = -1092
$ws.Range("N7").Value = -3402.75
$ws.Range("H122").Value = 3379.36
$ws.Range("I122").Value = 3861.75
$ws.Range("J122").Value = 2934.077
$ws.Range("K122").Value = 11585.25
$ws.Range("L122").Value = 8802.231
$ws.Range("M122").Value = -9135.25
$ws.Range("N122").Value = -13702.231
$ws.Range("H126").Value = 2783.8
$ws.Range("I126").Value = 1204
$ws.Range("J126").Value = 3178.75
$ws.Range("K126").Value = 3612
$ws.Range("L126").Value = 9536.25
$ws.Range("M126").Value = -1142
$ws.Range("N126").Value = -14476.25
$ws.Range("H128").Value = 30000
$ws.Range("J128").Value = 30000
$ws.Range("L128").Value = 30000
$ws.Range("N128").Value = -39960
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 40398.855
$ws.Range("J119").Value = 40398.855
$ws.Range("L119").Value = 40398.855
$ws.Range("N119").Value = -50074.855
$ws.Range("H124").Value = 39839.5
$ws.Range("J124").Value = 39839.5
$ws.Range("L124").Value = 39839.5
$ws.Range("N124").Value = -49659.5
